$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01751633333333333
$ws.Range("H2").Value = 0.052549
$ws.Range("I2").Value = 0.09281717406509865
$ws.Range("J2").Value = 0.09281717406509865
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 0.03550762788377777
$ws.Range("R2").Value = 0.319568650954
$ws.Range("S2").Value = 0.0006122484926913833
$ws.Range("T2").Value = 0.0006122484926913832
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01751633333333333
$ws.Range("H3").Value = 0.052549
$ws.Range("I3").Value = 0.09281717406509865
$ws.Range("J3").Value = 0.09281717406509865
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 4.491963933653777
$ws.Range("R3").Value = 40.427675402884
$ws.Range("S3").Value = 0.07745372787518857
$ws.Range("T3").Value = 0.07745372787518856
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01751633333333333
$ws.Range("H4").Value = 0.052549
$ws.Range("I4").Value = 0.09281717406509865
$ws.Range("J4").Value = 0.09281717406509865
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 0.855502373505889
$ws.Range("R4").Value = 7.699521361553001
$ws.Range("S4").Value = 0.0147511976972187
$ws.Range("T4").Value = 0.0147511976972187
$ws.Range("I5").Value = 0.4682437349423128
$ws.Range("J5").Value = 0.4682437349423127
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 0.1791287492504444
$ws.Range("R5").Value = 1.612158743254
$ws.Range("S5").Value = 0.003088668921653943
$ws.Range("T5").Value = 0.003088668921653942
$ws.Range("I6").Value = 0.4682437349423128
$ws.Range("J6").Value = 0.4682437349423127
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("S6").Value = 0.3907382786729456
$ws.Range("T6").Value = 0.3907382786729455
$ws.Range("I7").Value = 0.4682437349423128
$ws.Range("J7").Value = 0.4682437349423127
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 4.315835195989223
$ws.Range("R7").Value = 38.84251676390301
$ws.Range("S7").Value = 0.07441678734771322
$ws.Range("T7").Value = 0.07441678734771319
$ws.Range("G8").Value = 0.08283600000000001
$ws.Range("H8").Value = 0.248508
$ws.Range("I8").Value = 0.4389390909925887
$ws.Range("J8").Value = 0.4389390909925885
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 0.167918125752
$ws.Range("R8").Value = 1.511263131768
$ws.Range("S8").Value = 0.00289536715107329
$ws.Range("T8").Value = 0.002895367151073289
$ws.Range("G9").Value = 0.08283600000000001
$ws.Range("H9").Value = 0.248508
$ws.Range("I9").Value = 0.4389390909925887
$ws.Range("J9").Value = 0.4389390909925885
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 21.242820476592
$ws.Range("R9").Value = 191.185384289328
$ws.Range("S9").Value = 0.3662842491162033
$ws.Range("T9").Value = 0.3662842491162032
$ws.Range("G10").Value = 0.08283600000000001
$ws.Range("H10").Value = 0.248508
$ws.Range("I10").Value = 0.4389390909925887
$ws.Range("J10").Value = 0.4389390909925885
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 4.045732246764001
$ws.Range("R10").Value = 36.41159022087601
$ws.Range("S10").Value = 0.06975947472531212
$ws.Range("T10").Value = 0.06975947472531208
